$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 52.08999
$ws.Range("H2").Value = 156.26997
$ws.Range("I2").Value = 0.9401105828221099
$ws.Range("J2").Value = 0.9401105828221098
$ws.Range("M2").Value = 0.1994653333333334
$ws.Range("N2").Value = 0.598396
$ws.Range("O2").Value = 0.01676579960230272
$ws.Range("P2").Value = 0.01676579960230271
$ws.Range("Q2").Value = 10.39014721868
$ws.Range("R2").Value = 93.51132496812001
$ws.Range("S2").Value = 0.0157617056355995
$ws.Range("T2").Value = 0.0157617056355995
$ws.Range("G3").Value = 52.08999
$ws.Range("H3").Value = 156.26997
$ws.Range("I3").Value = 0.9401105828221099
$ws.Range("J3").Value = 0.9401105828221098
$ws.Range("O3").Value = 0.03203779682023726
$ws.Range("P3").Value = 0.03203779682023726
$ws.Range("Q3").Value = 19.85455113508
$ws.Range("R3").Value = 178.69096021572
$ws.Range("S3").Value = 0.03011907184100959
$ws.Range("T3").Value = 0.03011907184100959
$ws.Range("G4").Value = 52.08999
$ws.Range("H4").Value = 156.26997
$ws.Range("I4").Value = 0.9401105828221099
$ws.Range("J4").Value = 0.9401105828221098
$ws.Range("M4").Value = 0.2888043333333333
$ws.Range("N4").Value = 0.8664129999999999
$ws.Range("O4").Value = 0.02427507324719734
$ws.Range("P4").Value = 0.02427507324719734
$ws.Range("Q4").Value = 15.04381483529
$ws.Range("R4").Value = 135.39433351761
$ws.Range("S4").Value = 0.0228212532584721
$ws.Range("T4").Value = 0.02282125325847209
$ws.Range("G5").Value = 52.08999
$ws.Range("H5").Value = 156.26997
$ws.Range("I5").Value = 0.9401105828221099
$ws.Range("J5").Value = 0.9401105828221098
$ws.Range("M5").Value = 11.02772766666667
$ws.Range("N5").Value = 33.083183
$ws.Range("O5").Value = 0.9269213303302627
$ws.Range("P5").Value = 0.9269213303302626
$ws.Range("Q5").Value = 574.43422387939
$ws.Range("R5").Value = 5169.908014914509
$ws.Range("S5").Value = 0.8714085520870287
$ws.Range("T5").Value = 0.8714085520870285
$ws.Range("I6").Value = 0.009851545038079508
$ws.Range("J6").Value = 0.009851545038079508
$ws.Range("M6").Value = 0.1994653333333334
$ws.Range("N6").Value = 0.598396
$ws.Range("O6").Value = 0.01676579960230272
$ws.Range("P6").Value = 0.01676579960230271
$ws.Range("Q6").Value = 0.1088797479226667
$ws.Range("R6").Value = 0.979917731304
$ws.Range("S6").Value = 0.0001651690298815007
$ws.Range("T6").Value = 0.0001651690298815007
$ws.Range("I7").Value = 0.009851545038079508
$ws.Range("J7").Value = 0.009851545038079508
$ws.Range("O7").Value = 0.03203779682023726
$ws.Range("P7").Value = 0.03203779682023726
$ws.Range("S7").Value = 0.0003156217982954078
$ws.Range("T7").Value = 0.0003156217982954078
$ws.Range("I8").Value = 0.009851545038079508
$ws.Range("J8").Value = 0.009851545038079508
$ws.Range("M8").Value = 0.2888043333333333
$ws.Range("N8").Value = 0.8664129999999999
$ws.Range("O8").Value = 0.02427507324719734
$ws.Range("P8").Value = 0.02427507324719734
$ws.Range("Q8").Value = 0.1576461557846666
$ws.Range("R8").Value = 1.418815402062
$ws.Range("S8").Value = 0.0002391469773974436
$ws.Range("T8").Value = 0.0002391469773974435
$ws.Range("I9").Value = 0.009851545038079508
$ws.Range("J9").Value = 0.009851545038079508
$ws.Range("M9").Value = 11.02772766666667
$ws.Range("N9").Value = 33.083183
$ws.Range("O9").Value = 0.9269213303302627
$ws.Range("P9").Value = 0.9269213303302626
$ws.Range("Q9").Value = 6.019573368671332
$ws.Range("R9").Value = 54.17616031804199
$ws.Range("S9").Value = 0.009131607232505156
$ws.Range("T9").Value = 0.009131607232505156
$ws.Range("G10").Value = 1.744358333333333
$ws.Range("H10").Value = 5.233075
$ws.Range("I10").Value = 0.03148185917103467
$ws.Range("J10").Value = 0.03148185917103467
$ws.Range("M10").Value = 0.1994653333333334
$ws.Range("N10").Value = 0.598396
$ws.Range("O10").Value = 0.01676579960230272
$ws.Range("P10").Value = 0.01676579960230271
$ws.Range("Q10").Value = 0.3479390164111112
$ws.Range("R10").Value = 3.1314511477
$ws.Range("S10").Value = 0.0005278185419694833
$ws.Range("T10").Value = 0.0005278185419694831
$ws.Range("G11").Value = 1.744358333333333
$ws.Range("H11").Value = 5.233075
$ws.Range("I11").Value = 0.03148185917103467
$ws.Range("J11").Value = 0.03148185917103467
$ws.Range("O11").Value = 0.03203779682023726
$ws.Range("P11").Value = 0.03203779682023726
$ws.Range("Q11").Value = 0.6648772965222223
$ws.Range("R11").Value = 5.9838956687
$ws.Range("S11").Value = 0.001008609407644932
$ws.Range("T11").Value = 0.001008609407644932
$ws.Range("G12").Value = 1.744358333333333
$ws.Range("H12").Value = 5.233075
$ws.Range("I12").Value = 0.03148185917103467
$ws.Range("J12").Value = 0.03148185917103467
$ws.Range("M12").Value = 0.2888043333333333
$ws.Range("N12").Value = 0.8664129999999999
$ws.Range("O12").Value = 0.02427507324719734
$ws.Range("P12").Value = 0.02427507324719734
$ws.Range("Q12").Value = 0.5037782455527777
$ws.Range("R12").Value = 4.534004209975
$ws.Range("S12").Value = 0.0007642244373348181
$ws.Range("T12").Value = 0.0007642244373348178
$ws.Range("G13").Value = 1.744358333333333
$ws.Range("H13").Value = 5.233075
$ws.Range("I13").Value = 0.03148185917103467
$ws.Range("J13").Value = 0.03148185917103467
$ws.Range("M13").Value = 11.02772766666667
$ws.Range("N13").Value = 33.083183
$ws.Range("O13").Value = 0.9269213303302627
$ws.Range("P13").Value = 0.9269213303302626
$ws.Range("Q13").Value = 19.23630865308056
$ws.Range("R13").Value = 173.126777877725
$ws.Range("S13").Value = 0.02918120678408544
$ws.Range("T13").Value = 0.02918120678408543
$ws.Range("G14").Value = 1.028158333333334
$ws.Range("H14").Value = 3.084475
$ws.Range("I14").Value = 0.01855601296877595
$ws.Range("J14").Value = 0.01855601296877594
$ws.Range("M14").Value = 0.1994653333333334
$ws.Range("N14").Value = 0.598396
$ws.Range("O14").Value = 0.01676579960230272
$ws.Range("P14").Value = 0.01676579960230271
$ws.Range("Q14").Value = 0.2050819446777778
$ws.Range("R14").Value = 1.8457375021
$ws.Range("S14").Value = 0.0003111063948522278
$ws.Range("T14").Value = 0.0003111063948522276
$ws.Range("G15").Value = 1.028158333333334
$ws.Range("H15").Value = 3.084475
$ws.Range("I15").Value = 0.01855601296877595
$ws.Range("J15").Value = 0.01855601296877594
$ws.Range("O15").Value = 0.03203779682023726
$ws.Range("P15").Value = 0.03203779682023726
$ws.Range("Q15").Value = 0.3918914594555556
$ws.Range("R15").Value = 3.5270231351
$ws.Range("S15").Value = 0.0005944937732873314
$ws.Range("T15").Value = 0.0005944937732873313
$ws.Range("G16").Value = 1.028158333333334
$ws.Range("H16").Value = 3.084475
$ws.Range("I16").Value = 0.01855601296877595
$ws.Range("J16").Value = 0.01855601296877594
$ws.Range("M16").Value = 0.2888043333333333
$ws.Range("N16").Value = 0.8664129999999999
$ws.Range("O16").Value = 0.02427507324719734
$ws.Range("P16").Value = 0.02427507324719734
$ws.Range("Q16").Value = 0.2969365820194445
$ws.Range("R16").Value = 2.672429238175
$ws.Range("S16").Value = 0.0004504485739929798
$ws.Range("T16").Value = 0.0004504485739929797
$ws.Range("G17").Value = 1.028158333333334
$ws.Range("H17").Value = 3.084475
$ws.Range("I17").Value = 0.01855601296877595
$ws.Range("J17").Value = 0.01855601296877594
$ws.Range("M17").Value = 11.02772766666667
$ws.Range("N17").Value = 33.083183
$ws.Range("O17").Value = 0.9269213303302627
$ws.Range("P17").Value = 0.9269213303302626
$ws.Range("Q17").Value = 11.33825009821389
$ws.Range("R17").Value = 102.044250883925
$ws.Range("S17").Value = 0.01719996422664341
$ws.Range("T17").Value = 0.0171999642266434

Write-Host "Updated cells"